# Add a new "Coin" entry (row 11) with a hyperlinked source URL,
# mirroring the existing asset rows (e.g. row 4).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B11").Value = "Coin"
$ws.Range("D11").Value = "https://free3d.com/3d-model/coin-4532.html"

# Turn D11 into a real hyperlink pointing at the source page.
$ws.Hyperlinks.Add($ws.Range("D11"), "https://free3d.com/3d-model/coin-4532.html") | Out-Null

# Match the hyperlink styling already used elsewhere in the sheet (e.g. D4)
# instead of whatever new style Hyperlinks.Add auto-generated.
$ws.Range("D11").Style = $ws.Range("D4").Style

# Leave the selection where the author last left it.
$ws.Range("C17").Select() | Out-Null
